$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two data rows (originally "RM 232" at row 26 and "SC 92" at row 28) were removed
# entirely from the sheet. Delete bottom-up so row indices stay valid while deleting.
$ws.Rows("28").Delete()
$ws.Rows("26").Delete()

# After the row removal, several cells in column F (and a couple in column B)
# that were re-imputed / re-blanked differently than before need explicit updates.

# Newly filled-in (previously blank) cells
$ws.Range("F6").Value = 16.43
$ws.Range("F19").Value = 17.81
$ws.Range("F23").Value = 16.48
$ws.Range("B27").Value = -20.4
$ws.Range("F29").Value = 18.06

# Newly blanked (previously filled-in) cells
$ws.Range("F8").ClearContents()
$ws.Range("F21").ClearContents()
$ws.Range("B26").ClearContents()
$ws.Range("F27").ClearContents()
$ws.Range("B29").ClearContents()
